$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 2024 May
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 22

# Row 11: 2024 August
$ws.Range("B11").Value = 36
$ws.Range("C11").Value = 20

# Row 12: 2024 September
$ws.Range("B12").Value = 26
$ws.Range("C12").Value = 22

# Row 13: 2024 October
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 12
